$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The price/date records in rows 2-5 are cyclically shifted up by one row
# (row 3 -> row 2, row 4 -> row 3, row 5 -> row 4, row 2 -> row 5), while the
# market/category identifying columns (A,B,C,E,F,G,H,I,R) stay the same.

$dates = @(44691, 44687, 44221, 44692)
$volumen = @(100, 160, 250, 120)
$precioMin = @(3000, 3000, 1300, 3000)
$precioMax = @(3500, 3500, 1500, 3500)
$precioProm = @(3250, 3250, 1420, 3250)
$unidad = @("$/docena de matas", "$/docena de matas", "$/atado", "$/docena de matas")
$origen = @("Región Metropolitana", "Región Metropolitana", "Provincia de Diguillín", "Región Metropolitana")
$precioKg = @(542, 542, 1420, 542)
$kgUnidades = @(6, 6, 1, 6)

for ($i = 0; $i -lt 4; $i++) {
    $row = $i + 2
    $ws.Range("D$row").Value = $dates[$i]
    $ws.Range("J$row").Value = $volumen[$i]
    $ws.Range("K$row").Value = $precioMin[$i]
    $ws.Range("L$row").Value = $precioMax[$i]
    $ws.Range("M$row").Value = $precioProm[$i]
    $ws.Range("N$row").Value = $unidad[$i]
    $ws.Range("O$row").Value = $origen[$i]
    $ws.Range("P$row").Value = $precioKg[$i]
    $ws.Range("Q$row").Value = $kgUnidades[$i]
}
